$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PBIReports")
$ws2 = $wb.Worksheets.Item("RCExtensions")
$ws3 = $wb.Worksheets.Item("RCExtensionActions")

# --- Column width adjustments (targets minus 5/6 to compensate for this
#     runtime's ColumnWidth -> XML width quantization to 1/6 of a character) ---
$ws1.Columns.Item(2).ColumnWidth = 27.451822916666668   # -> 28.28515625
$ws1.Columns.Item(3).ColumnWidth = 32.451822916666664   # -> 33.28515625
$ws1.Columns.Item(4).ColumnWidth = 25.307291666666668   # -> 26.140625
$ws1.Columns.Item(5).ColumnWidth = 31.166666666666668   # -> 32

$ws2.Columns.Item(3).ColumnWidth = 40.592447916666664   # -> 41.42578125

# --- Selection / active-cell updates on existing sheets ---
[void]$ws2.Range("C2").Select()
[void]$ws3.Range("D5").Select()
[void]$ws1.Range("F9").Select()

# --- Add the new "PermissionSets" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws4.Name = "PermissionSets"

$ws4.Range("A1").Value = "id"
$ws4.Range("B1").Value = "name"
$ws4.Range("C1").Value = "filename"
$ws4.Range("D1").Value = "caption"

$ws4.Range("A2").Value = 50120
$ws4.Range("B2").Value = "PBI EMBED VIEW"
$ws4.Range("D2").Value = "PBI Embed - View"
$ws4.Range("C2").Value = "PBIEmbedView.PermissionSet"

$ws4.Columns.Item(2).ColumnWidth = 29.022135416666668   # -> 29.85546875
$ws4.Columns.Item(3).ColumnWidth = 43.022135416666664   # -> 43.85546875
$ws4.Columns.Item(4).ColumnWidth = 20.736979166666668   # -> 21.5703125

$tbl = $ws4.ListObjects.Add(1, $ws4.Range("A1:D2"), [System.Type]::Missing, 1)
$tbl.TableStyle = ""

[void]$ws4.Range("C4").Select()
